$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODAY")

$ws.Range("K7").Value = 136.838
$ws.Range("L7").Value = 103.29
$ws.Range("K8").Value = 209.948
$ws.Range("L8").Value = 252.22
